$d = $word.ActiveDocument

function Replace-ExactText($doc, $oldText, $newText) {
    # Locate the next occurrence of $oldText starting at the top of the
    # document each time (so positions recalculate after prior edits),
    # and replace only that run's characters - then "tickle" Bold on the
    # freshly written text (on -> off) so the paragraph's run list is
    # rebuilt/split instead of silently re-coalescing the edited run
    # into an adjacent run that happens to share identical formatting.
    $search = $doc.Content
    $search.Start = 0
    $search.End = $doc.Content.End
    if (-not $search.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
        throw "text not found: $oldText"
    }
    $startPos = $search.Start
    $target = $doc.Range($startPos, $search.End)
    $target.Text = $newText
    $endPos = $startPos + $newText.Length
    $written = $doc.Range($startPos, $endPos)
    $written.Font.Bold = 1
    $written.Font.Bold = 0
}

# Daily rate: 1.000 -> 5.000
Replace-ExactText $d "1.000" "5.000"

# Every remaining occurrence of the total amount: 3.000 -> 15.000
Replace-ExactText $d "3.000" "15.000"
Replace-ExactText $d "3.000" "15.000"
Replace-ExactText $d "3.000" "15.000"
Replace-ExactText $d "3.000" "15.000"
Replace-ExactText $d "3.000" "15.000"

Write-Output "applied kwitansi amount updates"
